# Site updated: 2021-01-21 12:51:03
#
# 1) Refresh the "datetimeFigureOut" date placeholder that lives on the
#    slide master and on every slide layout (PowerPoint re-stamps this
#    field text whenever the deck is saved) from 2021/1/19 -> 2021/1/21.
# 2) Resize/reposition the "PYLI" text box on slide 1 and bump its font
#    size from 180pt to 280pt.

$p = $ppt.ActivePresentation

$oldDate = "2021/1/19"
$newDate = "2021/1/21"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master date placeholder.
Update-DateShape $p.SlideMaster.Shapes

# Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DateShape $layouts.Item($L).Shapes
}

# PYLI text box on slide 1: move/resize and enlarge its text.
$s = $p.Slides.Item(1)
$pyli = $s.Shapes.Item(5)

$pyli.Left = 7017929 / 12700
$pyli.Top = 1228397 / 12700
$pyli.Width = 10851459 / 12700
$pyli.Height = 4401205 / 12700

$tr = $pyli.TextFrame.TextRange
$tr.Font.Size = 280
